$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row values (row 1) ---
$ws.Range("A1").Value = "Tätigkeit"
$ws.Range("B1").Value = "Beschreibung"
$ws.Range("C1").Value = "Dauer (in h)"
$ws.Range("D1").Value = "Ort des Dokuments (Falls vorhanden)"

# --- Border around the whole table (A1:D15), thin on all sides ---
$fullRange = $ws.Range("A1:D15")
$fullRange.Borders.LineStyle = 1

# --- Header formatting: white font on accent5 (blue) fill ---
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.ThemeColor = 2
$headerRange.Interior.ThemeColor = 9

# --- Column widths (autofit-like, matching content) ---
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 33.666666666666664

# --- Selection (column B selected, as in the saved view) ---
$ws.Range("B1:B1048576").Select()

Write-Host "done"
